$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Приходы")

# Лев (row 10): Акч (D10) 36 -> 38, * (E10) blank -> 2
$ws.Range("D10").Value = 38
$ws.Range("E10").Value = 2

# Настя ⊗ (row 12): Акч (D12) 47 -> 49
$ws.Range("D12").Value = 49

# Leave the cursor on A2, matching where editing left off
$ws.Range("A2").Select()

$wb.Save()
